{"js": "const body = context.document.body;\n\n// --- 1) Locate the unique \"on fir\" run, and the <del>/</del> runs that\n//        immediately wrap it (there are several <del> / </del> markers in\n//        the document, so we disambiguate using range position rather\n//        than index). -------------------------------------------------\nconst onFirResults = body.search(\"on fir\", { matchCase: true });\nonFirResults.load(\"items\");\n\nconst delOpenResults = body.search(\"<del>\", { matchCase: true });\ndelOpenResults.load(\"items\");\n\nconst delCloseResults = body.search(\"</del>\", { matchCase: true });\ndelCloseResults.load(\"items\");\n\nawait context.sync();\n\nif (onFirResults.items.length !== 1) {\n  throw new Error(\"expected exactly one match for 'on fir', found \" + onFirResults.items.length);\n}\n\nconst onFir = onFirResults.items[0];\nconst startOfOnFir = onFir.getRange(\"Start\");\nconst endOfOnFir = onFir.getRange(\"End\");\n\n// Find the <del> run whose end lands exactly at the start of \"on fir\".\nlet openCmps = [];\nfor (const item of delOpenResults.items) {\n  openCmps.push(item.getRange(\"End\").compareLocationWith(startOfOnFir));\n}\n// Find the </del> run whose start lands exactly at the end of \"on fir\".\nlet closeCmps = [];\nfor (const item of delCloseResults.items) {\n  closeCmps.push(item.getRange(\"Start\").compareLocationWith(endOfOnFir));\n}\nawait context.sync();\n\nlet delOpenRun = null;\nfor (let i = 0; i < openCmps.length; i++) {\n  if (openCmps[i].value === \"Equal\") {\n    delOpenRun = delOpenResults.items[i];\n    break;\n  }\n}\nlet delCloseRun = null;\nfor (let i = 0; i < closeCmps.length; i++) {\n  if (closeCmps[i].value === \"Equal\") {\n    delCloseRun = delCloseResults.items[i];\n    break;\n  }\n}\nif (!delOpenRun || !delCloseRun) {\n  throw new Error(\"could not locate the <del>/</del> pair surrounding 'on fir'\");\n}\n\n// --- 2) Apply the text-level edits, each confined to its own run so the\n//        existing run formatting (font/color/size) is preserved. -------\n\n// \"<del>\" -> \"<del><fr>\"\ndelOpenRun.insertText(\"<del><fr>\", Word.InsertLocation.replace);\n\n// \"on fir\" -> \"ou for\"\nonFir.insertText(\"ou for\", Word.InsertLocation.replace);\n\n// \"</del>\" -> \"</fr></del>\"\ndelCloseRun.insertText(\"</fr></del>\", Word.InsertLocation.replace);\n\n// \" such that the stem of tin, which is brittle &\" -> \" that the stem of tin, which is brittle &\"\nconst suchThatResults = body.search(\" such that the stem of tin, which is brittle &\", { matchCase: true });\nsuchThatResults.load(\"items\");\nawait context.sync();\nif (suchThatResults.items.length !== 1) {\n  throw new Error(\"expected exactly one match for the 'such that' run, found \" + suchThatResults.items.length);\n}\nsuchThatResults.items[0].insertText(\" that the stem of tin, which is brittle &\", Word.InsertLocation.replace);\n\n// \" delicate, is unable to support it\" -> \" delicate, would be unable to support\"\nconst delicateResults = body.search(\" delicate, is unable to support it\", { matchCase: true });\ndelicateResults.load(\"items\");\nawait context.sync();\nif (delicateResults.items.length !== 1) {\n  throw new Error(\"expected exactly one match for the 'delicate' run, found \" + delicateResults.items.length);\n}\ndelicateResults.items[0].insertText(\" delicate, would be unable to support\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1) Locate the unique \"on fir\" run and compute the character ranges\n#        of the immediately adjacent \"<del>\" / \"</del>\" markers. ------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"on fir\"\n$found = $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\nif (-not $found) {\n    throw \"could not find 'on fir'\"\n}\n$onFirRange = $d.Content\n$onFirRange.Start = $find.Parent.Start\n$onFirRange.End = $find.Parent.End\n\n$onFirStart = $onFirRange.Start\n$onFirEnd = $onFirRange.End\n\n# \"<del>\" immediately precedes \"on fir\"\n$delOpenRange = $d.Range($onFirStart - 5, $onFirStart)\nif ($delOpenRange.Text -ne \"<del>\") {\n    throw \"expected '<del>' immediately before 'on fir', got '$($delOpenRange.Text)'\"\n}\n\n# \"</del>\" immediately follows \"on fir\"\n$delCloseRange = $d.Range($onFirEnd, $onFirEnd + 6)\nif ($delCloseRange.Text -ne \"</del>\") {\n    throw \"expected '</del>' immediately after 'on fir', got '$($delCloseRange.Text)'\"\n}\n\n# --- 2) Apply the text-level edits, each confined to its own run so the\n#        existing run formatting (font/color/size) is preserved. -------\n\n# \"<del>\" -> \"<del><fr>\"\n$delOpenRange.Text = \"<del><fr>\"\n\n# \"on fir\" -> \"ou for\"   (re-fetch the range since prior edit shifted offsets)\n$onFirRange2 = $d.Range($onFirStart + 4, $onFirEnd + 4)\nif ($onFirRange2.Text -ne \"on fir\") {\n    throw \"expected 'on fir' after inserting '<fr>', got '$($onFirRange2.Text)'\"\n}\n$onFirRange2.Text = \"ou for\"\n\n# \"</del>\" -> \"</fr></del>\"   (re-fetch; offsets shifted by the two prior edits)\n$delCloseRange2 = $d.Range($onFirEnd + 4, $onFirEnd + 10)\nif ($delCloseRange2.Text -ne \"</del>\") {\n    throw \"expected '</del>', got '$($delCloseRange2.Text)'\"\n}\n$delCloseRange2.Text = \"</fr></del>\"\n\n# \" such that the stem of tin, which is brittle &\" -> \" that the stem of tin, which is brittle &\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \" such that the stem of tin, which is brittle &\"\n$found2 = $find2.Execute($find2.Text, $true, $false, $false, $false, $false, $true, 1, $false, \" that the stem of tin, which is brittle &\", 1)\nif (-not $found2) {\n    throw \"could not find/replace the 'such that' run\"\n}\n\n# \" delicate, is unable to support it\" -> \" delicate, would be unable to support\"\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Text = \" delicate, is unable to support it\"\n$found3 = $find3.Execute($find3.Text, $true, $false, $false, $false, $false, $true, 1, $false, \" delicate, would be unable to support\", 1)\nif (-not $found3) {\n    throw \"could not find/replace the 'delicate' run\"\n}\n\nWrite-Output \"done\"\n"}
